$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the transaction dates in column A (rows 2-26) forward, as in the
# updated upload of the workbook.
$ws.Range("A2").Value  = 45406
$ws.Range("A3").Value  = 45397
$ws.Range("A4").Value  = 45412
$ws.Range("A5").Value  = 45384
$ws.Range("A6").Value  = 45383
$ws.Range("A7").Value  = 45393
$ws.Range("A8").Value  = 45386
$ws.Range("A9").Value  = 45402
$ws.Range("A10").Value = 45403
$ws.Range("A11").Value = 45383
$ws.Range("A12").Value = 45405
$ws.Range("A13").Value = 45407
$ws.Range("A14").Value = 45394
$ws.Range("A15").Value = 45393
$ws.Range("A16").Value = 45409
$ws.Range("A17").Value = 45383
$ws.Range("A18").Value = 45384
$ws.Range("A19").Value = 45383
$ws.Range("A20").Value = 45410
$ws.Range("A21").Value = 45411
$ws.Range("A22").Value = 45397
$ws.Range("A23").Value = 45412
$ws.Range("A24").Value = 45399
$ws.Range("A25").Value = 45403
$ws.Range("A26").Value = 45387

# The re-uploaded workbook was saved with the cursor back on A1 instead of
# the previously selected E27.
$ws.Range("A1").Select()
